$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 400; this shifts the previous rows 400..477 down to 401..478
$ws.Rows.Item(400).Insert()

# Populate the newly inserted row 400 with the new record's data
$ws.Cells.Item(400, 1).Value = 11
$ws.Cells.Item(400, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(400, 3).Value = "Bíobío"
$ws.Cells.Item(400, 4).Value = 45211
$ws.Cells.Item(400, 5).Value = 8
$ws.Cells.Item(400, 6).Value = 100114013
$ws.Cells.Item(400, 7).Value = "Zanahoria"
$ws.Cells.Item(400, 8).Value = "Sin especificar"
$ws.Cells.Item(400, 9).Value = "Primera"
$ws.Cells.Item(400, 10).Value = 120
$ws.Cells.Item(400, 11).Value = 6000
$ws.Cells.Item(400, 12).Value = 6000
$ws.Cells.Item(400, 13).Value = 6000
$ws.Cells.Item(400, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(400, 15).Value = "Región Metropolitana"
$ws.Cells.Item(400, 16).Value = 300
$ws.Cells.Item(400, 17).Value = 20
$ws.Cells.Item(400, 18).Value = "Hortaliza"
